$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04241433333333333
$ws.Range("H2").Value = 0.127243
$ws.Range("I2").Value = 0.03157971673979624
$ws.Range("J2").Value = 0.04663324281038925
$ws.Range("M2").Value = 50.6072485
$ws.Range("N2").Value = 101.214497
$ws.Range("O2").Value = 0.2036920096625967
$ws.Range("P2").Value = 0.1535170070198019
$ws.Range("Q2").Value = 2.146472706961833
$ws.Range("R2").Value = 12.878836241771
$ws.Range("S2").Value = 0.006432535967304642
$ws.Range("T2").Value = 0.007158995863878653

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04241433333333333
$ws.Range("H3").Value = 0.127243
$ws.Range("I3").Value = 0.03157971673979624
$ws.Range("J3").Value = 0.04663324281038925
$ws.Range("M3").Value = 128.9086913333333
$ws.Range("N3").Value = 386.726074
$ws.Range("O3").Value = 0.5188519664463093
$ws.Range("P3").Value = 0.5865664620849566
$ws.Range("Q3").Value = 5.467576203775778
$ws.Range("R3").Value = 49.208185833982
$ws.Range("S3").Value = 0.01638519813026071
$ws.Range("T3").Value = 0.02735349625083876

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04241433333333333
$ws.Range("H4").Value = 0.127243
$ws.Range("I4").Value = 0.03157971673979624
$ws.Range("J4").Value = 0.04663324281038925
$ws.Range("M4").Value = 16.762851
$ws.Range("N4").Value = 50.288553
$ws.Range("O4").Value = 0.0674697579708304
$ws.Range("P4").Value = 0.07627512236628199
$ws.Range("Q4").Value = 0.710985149931
$ws.Range("R4").Value = 6.398866349378999
$ws.Range("S4").Value = 0.002130675845221433
$ws.Range("T4").Value = 0.00355695630169898

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04241433333333333
$ws.Range("H5").Value = 0.127243
$ws.Range("I5").Value = 0.03157971673979624
$ws.Range("J5").Value = 0.04663324281038925
$ws.Range("M5").Value = 35.4375075
$ws.Range("N5").Value = 70.87501499999999
$ws.Range("O5").Value = 0.1426344512705199
$ws.Range("P5").Value = 0.1074996220678108
$ws.Range("Q5").Value = 1.5030582556075
$ws.Range("R5").Value = 9.018349533644999
$ws.Range("S5").Value = 0.004504355568459289
$ws.Range("T5").Value = 0.005013055977913297

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04241433333333333
$ws.Range("H6").Value = 0.127243
$ws.Range("I6").Value = 0.03157971673979624
$ws.Range("J6").Value = 0.04663324281038925
$ws.Range("M6").Value = 3.403012
$ws.Range("N6").Value = 10.209036
$ws.Range("O6").Value = 0.01369697768069593
$ws.Range("P6").Value = 0.01548454715214769
$ws.Range("Q6").Value = 0.1443364853053334
$ws.Range("R6").Value = 1.299028367748
$ws.Range("S6").Value = 0.0004325466753476888
$ws.Range("T6").Value = 0.0007220946471550245

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.04241433333333333
$ws.Range("H7").Value = 0.127243
$ws.Range("I7").Value = 0.03157971673979624
$ws.Range("J7").Value = 0.04663324281038925
$ws.Range("M7").Value = 13.330536
$ws.Range("N7").Value = 39.991608
$ws.Range("O7").Value = 0.05365483696904789
$ws.Range("P7").Value = 0.06065723930900103
$ws.Range("Q7").Value = 0.5654057974159999
$ws.Range("R7").Value = 5.088652176744
$ws.Range("S7").Value = 0.00169440455320248
$ws.Range("T7").Value = 0.002828643768904533

# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.3006735
$ws.Range("H8").Value = 2.601347
$ws.Range("I8").Value = 0.9684202832602038
$ws.Range("J8").Value = 0.9533667571896107
$ws.Range("M8").Value = 50.6072485
$ws.Range("N8").Value = 101.214497
$ws.Range("O8").Value = 0.2036920096625967
$ws.Range("P8").Value = 0.1535170070198019
$ws.Range("Q8").Value = 65.82350703186475
$ws.Range("R8").Value = 263.294028127459
$ws.Range("S8").Value = 0.1972594736952921
$ws.Range("T8").Value = 0.1463580111559233

# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.3006735
$ws.Range("H9").Value = 2.601347
$ws.Range("I9").Value = 0.9684202832602038
$ws.Range("J9").Value = 0.9533667571896107
$ws.Range("M9").Value = 128.9086913333333
$ws.Range("N9").Value = 386.726074
$ws.Range("O9").Value = 0.5188519664463093
$ws.Range("P9").Value = 0.5865664620849566
$ws.Range("Q9").Value = 167.6681187369464
$ws.Range("R9").Value = 1006.008712421678
$ws.Range("S9").Value = 0.5024667683160485
$ws.Range("T9").Value = 0.5592129658341177

# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.3006735
$ws.Range("H10").Value = 2.601347
$ws.Range("I10").Value = 0.9684202832602038
$ws.Range("J10").Value = 0.9533667571896107
$ws.Range("M10").Value = 16.762851
$ws.Range("N10").Value = 50.288553
$ws.Range("O10").Value = 0.0674697579708304
$ws.Range("P10").Value = 0.07627512236628199
$ws.Range("Q10").Value = 21.8029960801485
$ws.Range("R10").Value = 130.817976480891
$ws.Range("S10").Value = 0.06533908212560896
$ws.Range("T10").Value = 0.072718166064583

# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.3006735
$ws.Range("H11").Value = 2.601347
$ws.Range("I11").Value = 0.9684202832602038
$ws.Range("J11").Value = 0.9533667571896107
$ws.Range("M11").Value = 35.4375075
$ws.Range("N11").Value = 70.87501499999999
$ws.Range("O11").Value = 0.1426344512705199
$ws.Range("P11").Value = 0.1074996220678108
$ws.Range("Q11").Value = 46.09262691130125
$ws.Range("R11").Value = 184.370507645205
$ws.Range("S11").Value = 0.1381300957020606
$ws.Range("T11").Value = 0.1024865660898975

# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.3006735
$ws.Range("H12").Value = 2.601347
$ws.Range("I12").Value = 0.9684202832602038
$ws.Range("J12").Value = 0.9533667571896107
$ws.Range("M12").Value = 3.403012
$ws.Range("N12").Value = 10.209036
$ws.Range("O12").Value = 0.01369697768069593
$ws.Range("P12").Value = 0.01548454715214769
$ws.Range("Q12").Value = 4.426207528582
$ws.Range("R12").Value = 26.557245171492
$ws.Range("S12").Value = 0.01326443100534824
$ws.Range("T12").Value = 0.01476245250499266

# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.3006735
$ws.Range("H13").Value = 2.601347
$ws.Range("I13").Value = 0.9684202832602038
$ws.Range("J13").Value = 0.9533667571896107
$ws.Range("M13").Value = 13.330536
$ws.Range("N13").Value = 39.991608
$ws.Range("O13").Value = 0.05365483696904789
$ws.Range("P13").Value = 0.06065723930900103
$ws.Range("Q13").Value = 17.338674915996
$ws.Range("R13").Value = 104.032049495976
$ws.Range("S13").Value = 0.05196043241584541
$ws.Range("T13").Value = 0.0578285955400965
